# Update "想去人数" (interest count) values in column F across the
# "展览" (Exhibition), "演出" (Performance) and "全部类型" (All types)
# sheets, per the refreshed scrape output.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 764
$wsExhibition.Range("F6").Value = 930
$wsExhibition.Range("F7").Value = 4503
$wsExhibition.Range("F8").Value = 347
$wsExhibition.Range("F9").Value = 495
$wsExhibition.Range("F10").Value = 837
$wsExhibition.Range("F17").Value = 1519
$wsExhibition.Range("F18").Value = 1398
$wsExhibition.Range("F19").Value = 602
$wsExhibition.Range("F22").Value = 214
$wsExhibition.Range("F24").Value = 88
$wsExhibition.Range("F25").Value = 1030
$wsExhibition.Range("F28").Value = 853
$wsExhibition.Range("F31").Value = 145
$wsExhibition.Range("F36").Value = 233
$wsExhibition.Range("F37").Value = 492

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F4").Value = 119
$wsShow.Range("F6").Value = 88

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 764
$wsAll.Range("F7").Value = 930
$wsAll.Range("F9").Value = 4505
$wsAll.Range("F10").Value = 347
$wsAll.Range("F11").Value = 496
$wsAll.Range("F13").Value = 119
$wsAll.Range("F14").Value = 837
$wsAll.Range("F17").Value = 88
$wsAll.Range("F24").Value = 1519
$wsAll.Range("F25").Value = 1398
$wsAll.Range("F26").Value = 602
$wsAll.Range("F29").Value = 214
$wsAll.Range("F32").Value = 88
$wsAll.Range("F33").Value = 1030
$wsAll.Range("F36").Value = 853
$wsAll.Range("F39").Value = 145
$wsAll.Range("F44").Value = 233
$wsAll.Range("F45").Value = 492
